$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (A12:D12) onto the new row (A13:D13)
# so the new row reuses the existing cell style instead of creating a new one.
$ws.Range("A12:D12").Copy()
$ws.Range("A13:D13").PasteSpecial(-4122)

# Fill in the new item's data
$ws.Range("A13").Value = "Problema ao cadastrar um radical para uma patente"
$ws.Range("B13").Value = "Defeito"
$ws.Range("C13").Value = "Em análise"

# Match the recorded selection after the edit
$ws.Range("A13:D13").Select()
